$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: MCH191-1
$ws.Range("A2").Value = "MCH191-1"
$ws.Range("C2").Value = "AA- AUSTRALIA, INFORMATION BULLETIN 1977-1992"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1977-1992"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 23P | GRAP COUNT NUMER: NONE"
$ws.Range("H2").Value = ""

# Row 3: MCH191-2
$ws.Range("A3").Value = "MCH191-2"
$ws.Range("C3").Value = "AA- AUSTRIA 1 BOOK, VARIOUS PAMPHLETS"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 23P | GRAP COUNT NUMER: NONE"
$ws.Range("H3").Value = ""

# Row 4: MCH191-3
$ws.Range("A4").Value = "MCH191-3"
$ws.Range("C4").Value = "WOMEN FOR JUSTICE IN SOUTHERN AFRICA- AMANDLA, ANTI- APARTHEID BEWEGENG, AAB NACHRICHTEN, ANNUAL REPORT TO FOUNDER, MEMORANDUM, LETTER TO DR. ODENDAAL, IMMSA, 1979 INTERNATIONAL YEAR OF THE CITY BY NELSON MANDELA, PRESENTATION FOR ROBBEN ISLAND RE-UNION, FOTOAVSSTELLUNG NAMIBIA: DAS AUSE DES ELECFANSTEIN"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1979"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "Series"
$ws.Range("F4").Value = "1 Box"
$ws.Range("G4").Value = "LOCATION: 24A | GRAP COUNT NUMER: NONE"
$ws.Range("H4").Value = ""

# Apply font formatting to the new data cells (skip column B which stays empty)
for ($r = 2; $r -le 4; $r++) {
    $ws.Range("A$r").Font.Name = "Calibri"
    $ws.Range("A$r").Font.Size = 10
    $ws.Range("A$r").Font.ThemeColor = 1

    $rng = $ws.Range("C$r`:H$r")
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 10
    $rng.Font.ThemeColor = 1
}

# Freeze the header row and select the data range, matching the saved view state
$ws.Range("A2:K4").Select()
$excel.ActiveWindow.FreezePanes = $true
